# GregShapiro_Resume.docx edits
# - Update cell phone number
# - Rewrite several CID bullet points
# - Rewrite NGP VAN "Led" bullet
# - Rewrite NGP VAN "Instructed on Git workflows" bullet

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Host "NOT FOUND: $old"
    }
    return $result
}

# 1. Phone number: (781) 354-7356 -> (617) 943-9870
Replace-Text "781) 354-7356 " "617) 943-9870 "

# 2. "Leading another..." bullet (CID, Product Architect)
Replace-Text `
    "Leading another geographically dispersed team of software engineers, UX, UI developers, visualization researcher" `
    "Leading a geographically dispersed team of software engineers, UX /UI designer, visualization researcher"

# 3. "Took over two software projects..." bullet (CID, Product Architect) -- full rewrite
Replace-Text `
    "Took over two software projects that were behind schedule and slipping fast, learned when the contractual deadlines were, and got the projects back on track, launching data driven, internationalized web apps for use by researchers policymakers in Mexico and Colombia" `
    "When I joined CID, I took over two software and data projects that were almost a year behind schedule and slipping fast. After learning what contractual deadlines and commitments were already in place, I worked with the team to focus and protect the most important aspects for creating a viable product. Between May and October, my team got the projects back on track, eventually launching data driven, visualization intensive, English/Spanish internationalized web apps for use by researchers and policymakers in Mexico and Colombia, along with launching smaller experimental visualizations along the way."

# 4. "Wrote the first shared Git workflow..." bullet (CID, Product Architect) -- full rewrite
Replace-Text `
    "Wrote the first shared Git workflow and instructions for software team members" `
    "While the team was using Git for most projects, I wrote up a standard workflow so that all developers could use the same branching, naming, and merging strategies"

# 5. "JIRA and github integration" -> capitalize Github
Replace-Text "JIRA and github integration" "JIRA and Github integration"

# 6. "where before" -> "whereas before"
Replace-Text `
    "Created single roadmap and cohesive planning for all software projects, where before the same team" `
    "Created single roadmap and cohesive planning for all software projects, whereas before the same team"

# 7. NGP VAN "Led a geographically dispersed team..." -> "Technical lead of a geographically dispersed team..."
Replace-Text "Led a geographically dispersed team of software engineers, " "Technical lead of a geographically dispersed team of software engineers, "

# 8. NGP VAN "Instructed on Git workflows, ensuring..." -> add "as we transitioned from TFS"
Replace-Text `
    "Instructed on Git workflows, ensuring kosher source control history " `
    "Instructed on Git workflows as we transitioned from TFS, ensuring kosher source control history "
